$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-like string (e.g. "61.40", "1.00")
# must be forced to Text format first, otherwise Excel auto-converts them to
# numbers and silently drops significant trailing/leading zeros (e.g. "61.40" -> 61.4).
$textCells = @(
    'D5',
    'D7',
    'D14',
    'D16',
    'D19',
    'D22',
    'D26',
    'D33',
    'D36',
    'D39',
    'D40',
    'D42',
    'D43',
    'D46',
    'D47',
    'D48',
    'D50',
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '38.666.15'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('D3').Value = '2.099.63'
$ws.Range('E3').Value = '  +3.59%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '229.47'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('E6').Value = '  +1.41%  '
$ws.Range('D7').Value = '61.40'
$ws.Range('E7').Value = '  +2.55%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').Value = '2.412.45'
$ws.Range('E12').Value = '  +3.66%  '
$ws.Range('E13').Value = '  +2.52%  '
$ws.Range('D14').Value = '22.49'
$ws.Range('E14').Value = '  +6.97%  '
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').Value = '5.48'
$ws.Range('E16').Value = '  +5.70%  '
$ws.Range('D17').Value = '2.072.68'
$ws.Range('E17').Value = '  +2.98%  '
$ws.Range('D18').Value = '38.557.65'
$ws.Range('E18').Value = '  +2.30%  '
$ws.Range('D19').Value = '70.82'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('E20').Value = '  +2.54%  '
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('D22').Value = '227.38'
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('E25').Value = '  +3.82%  '
$ws.Range('D26').Value = '170.10'
$ws.Range('E26').Value = '  +1.28%  '
$ws.Range('E27').Value = '  +1.08%  '
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('E29').Value = '  +1.70%  '
$ws.Range('E30').Value = '  +9.38%  '
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('E32').Value = '  +5.48%  '
$ws.Range('D33').Value = '4.76'
$ws.Range('E33').Value = '  +6.65%  '
$ws.Range('E34').Value = '  +2.31%  '
$ws.Range('E35').Value = '  +0.50%  '
$ws.Range('D36').Value = '6.48'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('E37').Value = '  +4.30%  '
$ws.Range('E38').Value = '  +4.13%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').Value = '18.41'
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('D41').Value = '1.541.60'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').Value = '100.07'
$ws.Range('E42').Value = '  +4.92%  '
$ws.Range('D43').Value = '0.0221'
$ws.Range('E43').Value = '  +2.74%  '
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').Value = '4.21'
$ws.Range('E46').Value = '  +3.34%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '7.57'
$ws.Range('E47').Value = '  +6.59%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = '1.12'
$ws.Range('E48').Value = '  +1.85%  '
$ws.Range('E49').Value = '  +3.77%  '
$ws.Range('D50').Value = '2.98'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').Value = '2.297.68'
$ws.Range('E51').Value = '  +3.64%  '

# Restore the text-forced cells back to the default (unstyled) look, now that the
# literal string value is locked in, so no stray number format lingers on them.
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
